# Implemented LeapProvider (converting, ...)
#
# Adds a new timelog entry (row 6) for 2014-05-31, 13:00-16:00, with the
# activity note "LeapProvider + Convert (untested), Contracts unfinished".
# The duration column (D) keeps its existing shared formula (=C-B) and the
# total in D31 (=SUM(D2:D30)) recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date (A6) - 2014-05-31, same day as row 5
$ws.Range("A6").Value = 41790
# from (B6) - 13:00
$ws.Range("B6").Value = 0.54166666666666663
# to (C6) - 16:00
$ws.Range("C6").Value = 0.66666666666666663
# activity (E6) - new note about the LeapProvider work
$ws.Range("E6").Value = "LeapProvider + Convert (untested), Contracts unfinished"

# Row 6's D cell already carries the shared formula "=C6-B6" (si="0");
# recalculate so it (and the D31 grand total) reflect the new duration.
$wb.Application.Calculate()
